$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18: new boundary row (TB_SUPCO2_DKISLBH_DKE_01 group, year 2010) ---
$ws.Range("D18").Value = 2010
$ws.Range("E18").Value = "UP"
$ws.Range("F18").Value = "CAP_BND"
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = "TB_SUPCO2_DKISLBH_DKE_01"

# --- Row 19 (year 2015) ---
$ws.Range("D19").Value = 2015
$ws.Range("E19").Value = "UP"
$ws.Range("F19").Value = "CAP_BND"
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Formula = "=M18"

# --- Row 20 (year 2020) ---
$ws.Range("D20").Value = 2020
$ws.Range("E20").Value = "UP"
$ws.Range("F20").Value = "CAP_BND"
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Formula = "=M19"

# --- Row 21 (year 2025) ---
$ws.Range("D21").Value = 2025
$ws.Range("E21").Value = "UP"
$ws.Range("F21").Value = "CAP_BND"
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Formula = "=M20"

# --- Row 22: new boundary row (TB_SUPCO2_DKISLBH_DKE_02 group, year 2010) ---
$ws.Range("D22").Value = 2010
$ws.Range("E22").Value = "UP"
$ws.Range("F22").Value = "CAP_BND"
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = "TB_SUPCO2_DKISLBH_DKE_02"

# --- Row 23 (year 2015) ---
$ws.Range("D23").Value = 2015
$ws.Range("E23").Value = "UP"
$ws.Range("F23").Value = "CAP_BND"
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Formula = "=M22"

# --- Row 24 (year 2020) ---
$ws.Range("D24").Value = 2020
$ws.Range("E24").Value = "UP"
$ws.Range("F24").Value = "CAP_BND"
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Formula = "=M23"

# --- Row 25 (year 2025) ---
$ws.Range("D25").Value = 2025
$ws.Range("E25").Value = "UP"
$ws.Range("F25").Value = "CAP_BND"
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Formula = "=M24"

$ws.Range("M23").Select()
